$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("criterio")
$ws.Columns.Item(6).Delete()
